$wb = $excel.ActiveWorkbook
$entities = $wb.Worksheets.Item("Entities")
$units = $wb.Worksheets.Item("Units")

# Row 3
$entities.Range("K3").Value = 0.055
$entities.Range("L3").Value = 18
$entities.Range("M3").Value = 18
$entities.Range("N3").Value = 0.019
$entities.Range("O3").Value = 0.016
$entities.Range("P3").Value = 24
$entities.Range("Q3").Value = 36
$entities.Range("R3").Value = 0.013000000000000001
$entities.Range("S3").Value = 0.0825
$entities.Range("T3").Value = 3.5
$entities.Range("V3").Value = 15
$entities.Range("Z3").Value = 39
$entities.Range("AA3").Value = 0.15
$entities.Range("AD3").Value = 0.06
$entities.Range("AE3").Value = 10
$entities.Range("AF3").Value = 10
$entities.Range("AG3").Value = 1
$entities.Range("AH3").Value = "Balloon"
$entities.Range("AJ3").Value = 0.013999999999999999
$entities.Range("AK3").Value = 15
$entities.Range("AL3").Value = 180
$entities.Range("AM3").Value = 1
$entities.Range("AN3").Value = 0.75
$entities.Range("AO3").Value = 0.3
$entities.Range("AP3").Value = 1.5
$entities.Range("AQ3").Value = 2
$entities.Range("AR3").Value = 0.025
$entities.Range("AS3").Value = 0.0075
$entities.Range("AT3").Value = 0.012
$entities.Range("AU3").Value = 6
$entities.Range("AV3").Value = 6
$entities.Range("AW3").Value = 0.0025
$entities.Range("AX3").Value = 0.012
$entities.Range("AY3").Value = 6
$entities.Range("AZ3").Value = 6
$entities.Range("BA3").Value = 0.0025
$entities.Range("BB3").Value = 0.0005
$entities.Range("BF3").Value = 1.6

# Row 5
$entities.Range("D5").ClearContents()

# Row 6
$entities.Range("K6").Value = 0.06
$entities.Range("L6").Value = 10
$entities.Range("M6").Value = 7
$entities.Range("N6").Value = 0.0148
$entities.Range("O6").Value = 0.0443
$entities.Range("P6").Value = 10
$entities.Range("Q6").Value = 10
$entities.Range("R6").Value = 0.024300000000000002
$entities.Range("T6").Value = 4.86
$entities.Range("AD6").Value = 0.2335
$entities.Range("AK6").Value = 9
$entities.Range("AN6").Value = 0.6751999999999999
$entities.Range("AP6").Value = 2.2
$entities.Range("AQ6").Value = 2
$entities.Range("AR6").Value = 0.0345
$entities.Range("AS6").Value = 0.035
$entities.Range("AT6").Value = 0.0085
$entities.Range("AW6").Value = 0.001
$entities.Range("AX6").Value = 0.0085
$entities.Range("BA6").Value = 0.001

# Row 10
$entities.Range("O10").Value = 0.016
$entities.Range("T10").Value = 3.5
$entities.Range("Y10").Value = 0
$entities.Range("Z10").Value = 39
$entities.Range("AE10").Value = 10
$entities.Range("AF10").Value = 10
$entities.Range("AJ10").Value = 0.013999999999999999
$entities.Range("AK10").Value = 15
$entities.Range("AO10").Value = 0.3
$entities.Range("AP10").Value = 1.5
$entities.Range("AQ10").Value = 2
$entities.Range("AT10").Value = 0.012
$entities.Range("AX10").Value = 0.012
$entities.Range("BB10").Value = 0.0005

# Row 11
$entities.Range("AJ11").Value = 0.025
$entities.Range("AO11").Value = 0.5
$entities.Range("BB11").Value = 0.0005

# Row 12
$entities.Range("BB12").Value = 0.0005

# Row 13
$entities.Range("O13").Value = 0.016
$entities.Range("R13").Value = 0.013000000000000001
$entities.Range("S13").Value = 0.06849999999999999
$entities.Range("Y13").Value = 0
$entities.Range("AJ13").Value = 0.0125
$entities.Range("AK13").Value = 15
$entities.Range("AO13").Value = 0.35
$entities.Range("AQ13").Value = 2.5
$entities.Range("BB13").Value = 0.0005
$entities.Range("BC13").Value = 0
$entities.Range("BD13").Value = 0

# Row 14
$entities.Range("O14").Value = 0.018000000000000002
$entities.Range("R14").Value = 0.015
$entities.Range("S14").Value = 0.0655
$entities.Range("T14").Value = 2
$entities.Range("AJ14").Value = 0.022000000000000002
$entities.Range("AO14").Value = 0.45
$entities.Range("AP14").Value = 0.75
$entities.Range("BB14").Value = 0.0005


# Add new Unit 2 / OGANITA lease row to Units sheet
$units.Cells.Item(34, 2).Value = 32
$units.Cells.Item(34, 3).Value = "Unit 2"
$units.Cells.Item(34, 4).Value = 9
$units.Cells.Item(34, 5).Value = 1
$units.Cells.Item(34, 6).Value = 2017
$units.Cells.Item(34, 7).Value = 12
$units.Cells.Item(34, 8).Value = 2050
$units.Cells.Item(34, 9).Value = "Occupied"
$units.Cells.Item(34, 10).Value = "OGANITA"
$units.Cells.Item(34, 11).Value = 4500000
$units.Cells.Item(34, 22).Value = 250000
$units.Cells.Item(34, 23).Value = 18

# Restore view state: scroll back to left, select G16 on Entities sheet
$entities.Activate()
$entities.Range("G16").Select()
